# Update "想去人数" (interested-count) figures for two events across two sheets,
# reflecting refreshed numbers from the latest data pull.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions): rows 2-4 hold F (想去人数) values that changed.
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 412
$wsExhibit.Range("F3").Value = 2430
$wsExhibit.Range("F4").Value = 115

# Sheet "全部类型" (All types): same three events appear again at rows 2, 7, 8.
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 412
$wsAll.Range("F7").Value = 2430
$wsAll.Range("F8").Value = 115
